$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 0.002
$ws.Range("K4").Value = 1016
$ws.Range("L4").Value = 0.002032
